# Daily attendance processing - 2025-11-28 08:59:32
# Rotate the "Recorded By" (column G) value's comma-separated list left by
# one position (move the first entry to the end) for every data row where
# the list either already ends with an exact "System" entry, or contains a
# lowercase "system" entry anywhere in the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2

    if ($current -eq $null) {
        continue
    }

    $parts = $current -split ", "
    if ($parts.Count -lt 2) {
        continue
    }

    $shouldRotate = $false
    $lastPart = $parts[$parts.Count - 1]
    if ($lastPart.Equals("System")) {
        $shouldRotate = $true
    }
    foreach ($p in $parts) {
        if ($p.Equals("system")) {
            $shouldRotate = $true
        }
    }

    if ($shouldRotate) {
        $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
        $cell.Value = ($rotated -join ", ")
    }
}
